# Add two new order line items to the "orders" sheet, matching the
# checkout flow now writing each purchased item (with its own order id)
# into its own row, and set a narrower, explicit width for columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: first item of the order
$ws.Range("A2").Value = "d5531907-1cc6-4abb-8302-5182051301be"
$ws.Range("B2").Value = "HERO-2020 HOODIES"
$ws.Range("C2").Value = "M"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 36.0099983215332

# Row 3: second item tied to the same order/product, distinct order id
$ws.Range("A3").Value = "df4b0000-1a5e-4d9a-af8b-73db6ecf3a59"
$ws.Range("B3").Value = "HERO-2020 HOODIES"
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 36.0099983215332

# Explicit column width for A:E
$ws.Columns("A:E").ColumnWidth = 7.5
